$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 2 new rows for the new worker (EUDIN HERNANDEZ MORALES) above the
#     existing table body, pushing the rest of the rows (and the footer) down
#     by 2. Only touch columns B:J so we don't balloon the row out to XFD. ---
$ws.Range("B16:J17").Insert()

# Copy the formatting (borders, number formats, etc.) of an ordinary data row
# (row 18, a "PE" row that kept the plain interior style) onto the two blank
# rows we just inserted.
$ws.Range("B18:J18").Copy()
$ws.Range("B16:J17").PasteSpecial(-4122)

# --- New worker rows ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "12602950"
$ws.Range("D16").Value = "EUDIN HERNANDEZ MORALES"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 35112
$ws.Range("G16").Value = 877803

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "12602950"
$ws.Range("D17").Value = "EUDIN HERNANDEZ MORALES"
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 35112
$ws.Range("G17").Value = 877803

# --- ABRAHAM JOSE GIL HERNANDEZ block: periods reordered (desc) ---
$ws.Range("E18").Value = "2102"
$ws.Range("F18").Value = 35112

$ws.Range("E19").Value = "2101"

$ws.Range("E20").Value = "2012"

$ws.Range("E21").Value = "2011"
$ws.Range("F21").Value = 7022

# --- JOSE LUIS CARRASCAL MACHADO block: periods reordered (desc) ---
$ws.Range("E22").Value = "2107"
$ws.Range("F22").Value = 19382

$ws.Range("E23").Value = "2106"

$ws.Range("E25").Value = "2104"

$ws.Range("E26").Value = "2103"
$ws.Range("F26").Value = 36341

# --- Summary cells ---
$ws.Range("E11").Value = 347328
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 11
